$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3.3333333  # was 2.75
$ws.Range("I11").Value = 3.3333333  # was 2.75
$ws.Range("K11").Value = 3.3333333  # was 2.75
$ws.Range("M11").Value = 136.6666667  # was 137.25
$ws.Range("H15").Value = 1453.5319  # was 1457.0426
$ws.Range("I15").Value = 1453.5319  # was 1457.0426
$ws.Range("K15").Value = 4360.5957  # was 4371.1278
$ws.Range("M15").Value = -4191.5957  # was -4202.1278
$ws.Range("H40").Value = 5428.0713  # was 5266.2
$ws.Range("H52").Value = 465  # was 467
$ws.Range("I52").Value = 465  # was 467
$ws.Range("K52").Value = 1395  # was 1401
$ws.Range("M52").Value = -1235  # was -1241
$ws.Range("H70").Value = 18706.25  # was 16838.889
$ws.Range("I70").Value = 1971  # was 1947.3334
$ws.Range("K70").Value = 5913  # was 5842.0002
$ws.Range("M70").Value = -5643  # was -5572.0002
$ws.Range("H73").Value = 18706.25  # was 16838.889
$ws.Range("I73").Value = 1971  # was 1947.3334
$ws.Range("K73").Value = 5913  # was 5842.0002
$ws.Range("M73").Value = -4977  # was -4906.0002
$ws.Range("H135").Value = 52633668  # was 58825772
$ws.Range("I135").Value = 71430744  # was 83335736
$ws.Range("K135").Value = 642876696  # was 750021624
$ws.Range("M135").Value = -642874161  # was -750019089

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4128.25  # was 4243.38
$ws.Range("I32").Value = 1616.829  # was 1768.3158
$ws.Range("K32").Value = 1616.829  # was 1768.3158
$ws.Range("M32").Value = -1329.829  # was -1481.3158
$ws.Range("H45").Value = 6240.7617  # was 6319.5264
$ws.Range("I45").Value = 4278.3335  # was 4370.2144
$ws.Range("J45").Value = 11146.833  # was 11777.6
$ws.Range("K45").Value = 4278.3335  # was 4370.2144
$ws.Range("L45").Value = 11146.833  # was 11777.6
$ws.Range("M45").Value = -3901.3335  # was -3993.2144
$ws.Range("N45").Value = -11900.833  # was -12531.6
$ws.Range("H61").Value = 5291.5854  # was 5758.5674
$ws.Range("I61").Value = 2232  # was 2425.8462
$ws.Range("K61").Value = 2232  # was 2425.8462
$ws.Range("M61").Value = -2020  # was -2213.8462
$ws.Range("H63").Value = 10500.5  # was 7999.6924
$ws.Range("I63").Value = 4502  # was 4166
$ws.Range("J63").Value = 12500  # was 11285.714
$ws.Range("K63").Value = 4502  # was 4166
$ws.Range("L63").Value = 12500  # was 11285.714
$ws.Range("M63").Value = -3816  # was -3480
$ws.Range("N63").Value = -13872  # was -12657.714
$ws.Range("H66").Value = 10500.5  # was 7999.6924
$ws.Range("I66").Value = 4502  # was 4166
$ws.Range("J66").Value = 12500  # was 11285.714
$ws.Range("K66").Value = 22510  # was 20830
$ws.Range("L66").Value = 62500  # was 56428.57
$ws.Range("M66").Value = -19078  # was -17398
$ws.Range("N66").Value = -69364  # was -63292.57
$ws.Range("H74").Value = 2597.1765  # was 2627
$ws.Range("J74").Value = 2864.6667  # was 2951.8076
$ws.Range("L74").Value = 2864.6667  # was 2951.8076
$ws.Range("N74").Value = -4612.6667  # was -4699.8076
$ws.Range("H77").Value = 2597.1765  # was 2627
$ws.Range("J77").Value = 2864.6667  # was 2951.8076
$ws.Range("L77").Value = 14323.3335  # was 14759.038
$ws.Range("N77").Value = -23059.3335  # was -23495.038
$ws.Range("H102").Value = 1940.1613  # was 2065.5862
$ws.Range("I102").Value = 1453.6207  # was 1552.2963
$ws.Range("K102").Value = 1453.6207  # was 1552.2963
$ws.Range("M102").Value = 168.3793000000001  # was 69.70370000000003
$ws.Range("H119").Value = 179445.75  # was 200080.86
$ws.Range("J119").Value = 179445.75  # was 200080.86
$ws.Range("L119").Value = 179445.75  # was 200080.86
$ws.Range("N119").Value = -189121.75  # was -209756.86
$ws.Range("H122").Value = 3942.8823  # was 4064.375
$ws.Range("I122").Value = 2338.1667  # was 2406
$ws.Range("K122").Value = 7014.500100000001  # was 7218
$ws.Range("M122").Value = -4564.500100000001  # was -4768
$ws.Range("H135").Value = 79998.336  # was 79999
$ws.Range("J135").Value = 79998.336  # was 79999
$ws.Range("L135").Value = 79998.336  # was 79999
$ws.Range("N135").Value = -90138.336  # was -90139
$ws.Range("H136").Value = 5291.5854  # was 5758.5674
$ws.Range("I136").Value = 2232  # was 2425.8462
$ws.Range("K136").Value = 6696  # was 7277.5386
$ws.Range("M136").Value = -4146  # was -4727.5386

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2128.3809  # was 1853.9259
$ws.Range("I105").Value = 2087.8333  # was 1789.2084
$ws.Range("K105").Value = 2087.8333  # was 1789.2084
$ws.Range("M105").Value = -340.8332999999998  # was -42.20839999999998
$ws.Range("H135").Value = 99987.5  # was 99988
$ws.Range("J135").Value = 99987.5  # was 99988
$ws.Range("L135").Value = 99987.5  # was 99988
$ws.Range("N135").Value = -110127.5  # was -110128

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1989.3529  # was 2108
$ws.Range("I58").Value = 1393.3334  # was 1438.1818
$ws.Range("J58").Value = 3419.8  # was 3950
$ws.Range("K58").Value = 1393.3334  # was 1438.1818
$ws.Range("L58").Value = 3419.8  # was 3950
$ws.Range("M58").Value = -1190.3334  # was -1235.1818
$ws.Range("N58").Value = -3825.8  # was -4356
$ws.Range("H132").Value = 1473.1945  # was 1492.4286
$ws.Range("I132").Value = 1488.2142  # was 1513.7037
$ws.Range("K132").Value = 4464.642599999999  # was 4541.1111
$ws.Range("M132").Value = -1934.642599999999  # was -2011.1111
$ws.Range("H136").Value = 1989.3529  # was 2108
$ws.Range("I136").Value = 1393.3334  # was 1438.1818
$ws.Range("J136").Value = 3419.8  # was 3950
$ws.Range("K136").Value = 4180.0002  # was 4314.5454
$ws.Range("L136").Value = 10259.4  # was 11850
$ws.Range("M136").Value = -1630.0002  # was -1764.5454
$ws.Range("N136").Value = -15359.4  # was -16950

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 50105336  # was 41157956
$ws.Range("I4").Value = 56321136  # was 45056908
$ws.Range("K4").Value = 168963408  # was 135170724
$ws.Range("M4").Value = -168963296  # was -135170612
$ws.Range("H131").Value = 11309.083  # was 11370.917
$ws.Range("I131").Value = 466.33334  # was 475
$ws.Range("J131").Value = 12858.048  # was 12361.454
$ws.Range("K131").Value = 1399.00002  # was 1425
$ws.Range("L131").Value = 38574.144  # was 37084.362
$ws.Range("M131").Value = 3640.99998  # was 3615
$ws.Range("N131").Value = -48654.144  # was -47164.362

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 472.42856  # was 476
$ws.Range("I107").Value = 381.75  # was 450.66666
$ws.Range("J107").Value = 593.3333  # was 495
$ws.Range("K107").Value = 381.75  # was 450.66666
$ws.Range("L107").Value = 593.3333  # was 495
$ws.Range("M107").Value = 1538.25  # was 1469.33334
$ws.Range("N107").Value = -4433.3333  # was -4335
$ws.Range("H113").Value = 2507.3076  # was 2699.3635
$ws.Range("I113").Value = 1465.9166  # was 1580
$ws.Range("J113").Value = 3399.9285  # was 3632.1667
$ws.Range("K113").Value = 1465.9166  # was 1580
$ws.Range("L113").Value = 3399.9285  # was 3632.1667
$ws.Range("M113").Value = 704.0834  # was 590
$ws.Range("N113").Value = -7739.9285  # was -7972.1667
$ws.Range("H132").Value = 635857.7  # was 671127.7
$ws.Range("I132").Value = 671016.4399999999  # was 710429.3
$ws.Range("K132").Value = 2013049.32  # was 2131287.9
$ws.Range("M132").Value = -2010519.32  # was -2128757.9

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3523.4375  # was 3498.5
$ws.Range("J46").Value = 3540.3225  # was 3514.5806
$ws.Range("L46").Value = 3540.3225  # was 3514.5806
$ws.Range("N46").Value = -3916.3225  # was -3890.5806
$ws.Range("H82").Value = 2286.5715  # was 1971.125
$ws.Range("I82").Value = 1979.7778  # was 1980
$ws.Range("J82").Value = 2838.8  # was 1959.7142
$ws.Range("K82").Value = 1979.7778  # was 1980
$ws.Range("L82").Value = 2838.8  # was 1959.7142
$ws.Range("M82").Value = -1618.7778  # was -1619
$ws.Range("N82").Value = -3560.8  # was -2681.7142
$ws.Range("H85").Value = 2286.5715  # was 1971.125
$ws.Range("I85").Value = 1979.7778  # was 1980
$ws.Range("J85").Value = 2838.8  # was 1959.7142
$ws.Range("K85").Value = 1979.7778  # was 1980
$ws.Range("L85").Value = 2838.8  # was 1959.7142
$ws.Range("M85").Value = -731.7778000000001  # was -732
$ws.Range("N85").Value = -5334.8  # was -4455.7142
$ws.Range("H127").Value = 139796.28  # was 133252.44
$ws.Range("J127").Value = 139796.28  # was 133252.44
$ws.Range("L127").Value = 139796.28  # was 133252.44
$ws.Range("N127").Value = -149716.28  # was -143172.44
$ws.Range("H132").Value = 1332969.9  # was 1649764
$ws.Range("I132").Value = 1731285.9  # was 2163502.8
$ws.Range("J132").Value = 5250  # was 5800
$ws.Range("K132").Value = 5193857.699999999  # was 6490508.399999999
$ws.Range("L132").Value = 15750  # was 17400
$ws.Range("M132").Value = -5191327.699999999  # was -6487978.399999999
$ws.Range("N132").Value = -20810  # was -22460

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 859.9  # was 933
$ws.Range("I81").Value = 844.3333  # was 933
$ws.Range("J81").Value = 1000  # was 0
$ws.Range("K81").Value = 1688.6666  # was 1866
$ws.Range("L81").Value = 2000  # was 0
$ws.Range("M81").Value = -627.6666  # was -805
$ws.Range("N81").Value = -4122  # new cell
$ws.Range("H84").Value = 859.9  # was 933
$ws.Range("I84").Value = 844.3333  # was 933
$ws.Range("J84").Value = 1000  # was 0
$ws.Range("K84").Value = 8443.333000000001  # was 9330
$ws.Range("L84").Value = 10000  # was 0
$ws.Range("M84").Value = -3139.333000000001  # was -4026
$ws.Range("N84").Value = -20608  # new cell
$ws.Range("H96").Value = 1807.909  # was 2068.2856
$ws.Range("I96").Value = 1808.8889  # was 2112.3333
$ws.Range("J96").Value = 1803.5  # was 1804
$ws.Range("K96").Value = 1808.8889  # was 2112.3333
$ws.Range("L96").Value = 1803.5  # was 1804
$ws.Range("M96").Value = -435.8888999999999  # was -739.3332999999998
$ws.Range("N96").Value = -4549.5  # was -4550
$ws.Range("H107").Value = 2466.3572  # was 2076.1177
$ws.Range("I107").Value = 1231.1818  # was 1022
$ws.Range("K107").Value = 3693.5454  # was 3066
$ws.Range("M107").Value = -1773.5454  # was -1146
$ws.Range("H132").Value = 1506430  # was 1649797.2
$ws.Range("I132").Value = 2660514.8  # was 3144049.5
$ws.Range("K132").Value = 7981544.399999999  # was 9432148.5
$ws.Range("M132").Value = -7979014.399999999  # was -9429618.5
$ws.Range("H135").Value = 0  # was 51000
$ws.Range("J135").Value = 0  # was 51000
$ws.Range("L135").Value = 0  # was 51000
$ws.Range("N135").ClearContents()  # remove cell (was -61140)
